# Applies the "cryptos list" price/volume refresh described by the commit
# "Updated cryptos list on Wed Nov 29 15:15:52 UTC 2023 with GitHub Actions".
#
# Column D ("Price") values such as "226.77" or "0.999" look like plain numbers,
# and Excel would silently convert them to numeric cells (dropping formatting like
# trailing zeros and introducing binary floating point noise). To keep them as the
# literal text strings found in the sheet, those values are written with a leading
# apostrophe (forces text entry) and the cell style is immediately reset back to
# "Normal" so no stray quote-prefix formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "37.711.08"
$ws.Range("E2").Value = "  +0.78%  "
# Row 3
$ws.Range("D3").Value = "2.027.92"
$ws.Range("E3").Value = "  -0.24%  "
# Row 4
$ws.Range("E4").Value = "  +0.37%  "
# Row 5
$ws.Range("D5").Value = "'226.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.22%  "
# Row 6
$ws.Range("D6").Value = "'0.611"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.45%  "
# Row 7
$ws.Range("D7").Value = "'59.86"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +6.76%  "
# Row 8
$ws.Range("E8").Value = "  +0.08%  "
# Row 9
$ws.Range("D9").Value = "'0.379"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.59%  "
# Row 10
$ws.Range("D10").Value = "'0.0808"
$ws.Range("D10").Style = "Normal"
# Row 11
$ws.Range("D11").Value = "'0.103"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.47%  "
# Row 12
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "2.329.57"
$ws.Range("E12").Value = "  -0.23%  "
# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'14.52"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.81%  "
# Row 14
$ws.Range("D14").Value = "'20.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.92%  "
# Row 15
$ws.Range("D15").Value = "'0.751"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.03%  "
# Row 16
$ws.Range("D16").Value = "'5.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.06%  "
# Row 17
$ws.Range("D17").Value = "2.032.44"
$ws.Range("E17").Value = "  -0.29%  "
# Row 18
$ws.Range("D18").Value = "37.706.30"
$ws.Range("E18").Value = "  +0.84%  "
# Row 19
$ws.Range("D19").Value = "'6.07"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.34%  "
# Row 20
$ws.Range("D20").Value = "'69.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.61%  "
# Row 21
$ws.Range("D21").Value = "0.0₃0820"
$ws.Range("E21").Value = "  -0.59%  "
# Row 22
$ws.Range("D22").Value = "'222.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.32%  "
# Row 23
$ws.Range("E23").Value = "  +0.15%  "
# Row 24
$ws.Range("E24").Value = "  -0.74%  "
# Row 25
$ws.Range("D25").Value = "'2.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.00%  "
# Row 26
$ws.Range("D26").Value = "'165.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.47%  "
# Row 27
$ws.Range("D27").Value = "'9.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.52%  "
# Row 28
$ws.Range("D28").Value = "'0.129"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.12%  "
# Row 29
$ws.Range("D29").Value = "'18.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.03%  "
# Row 30
$ws.Range("D30").Value = "'1.27"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.58%  "
# Row 31
$ws.Range("D31").Value = "'0.119"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.41%  "
# Row 32
$ws.Range("D32").Value = "'4.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.50%  "
# Row 33
$ws.Range("E33").Value = "  +2.52%  "
# Row 34
$ws.Range("D34").Value = "'4.49"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.30%  "
# Row 35
$ws.Range("D35").Value = "'0.0601"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.89%  "
# Row 36
$ws.Range("D36").Value = "'6.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.48%  "
# Row 37
$ws.Range("D37").Value = "'2.26"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.69%  "
# Row 38
$ws.Range("D38").Value = "'3.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.29%  "
# Row 39
$ws.Range("E39").Value = "  -0.06%  "
# Row 40
$ws.Range("D40").Value = "1.527.45"
$ws.Range("E40").Value = "  +3.58%  "
# Row 41
$ws.Range("D41").Value = "'0.0217"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.85%  "
# Row 42
$ws.Range("D42").Value = "'96.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.61%  "
# Row 43
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").Value = "'16.45"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.76%  "
# Row 44
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").Value = "'2.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.54%  "
# Row 45
$ws.Range("D45").Value = "'0.0916"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.55%  "
# Row 46
$ws.Range("E46").Value = "  -0.31%  "
# Row 47
$ws.Range("D47").Value = "'4.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.20%  "
# Row 48
$ws.Range("D48").Value = "'2.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.59%  "
# Row 49
$ws.Range("D49").Value = "'0.999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.06%  "
# Row 50
$ws.Range("E50").Value = "  -1.18%  "
# Row 51
$ws.Range("D51").Value = "2.224.01"
$ws.Range("E51").Value = "  +0.08%  "
